# Add a new, empty paragraph (no runs) right after the last empty
# "spacer" paragraph that precedes the "1. Identificación de las
# necesidades" heading. The new paragraph keeps the same spacing/font
# as its neighbour but marks its (empty) paragraph mark as underlined
# (<w:u w:val="single"/> in w:pPr/w:rPr), matching the target diff:
#
#   <w:p>
#     <w:pPr>
#       <w:spacing w:before="320" w:after="320" w:line="360" w:lineRule="auto"/>
#       <w:rPr>
#         <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
#         <w:sz w:val="24"/>
#         <w:szCs w:val="24"/>
#         <w:u w:val="single"/>
#       </w:rPr>
#     </w:pPr>
#   </w:p>

$d = $word.ActiveDocument

# Locate the heading that follows the target blank paragraph. The TOC
# entry for this same text lives inside a TOC field and is not matched
# by Find, so this reliably lands on the real heading occurrence.
$needle = "1. Identificaci"
$rng = $d.Content
$found = $rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 0, $false, "", 0)

$insertionPoint = -1

if ($found) {
    $insertionPoint = $rng.Start
} else {
    # Fallback: scan paragraphs for the heading text and use the start
    # of its range as the insertion point.
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "1. Identificaci*") {
            $insertionPoint = $p.Range.Start
            break
        }
    }
}

if ($insertionPoint -ge 0) {
    # Collapsed (zero-length) range positioned right before the heading,
    # i.e. right after the preceding blank paragraph's own paragraph
    # mark - inserting here adds a new paragraph between the two
    # without disturbing either neighbour.
    $target = $d.Range($insertionPoint, $insertionPoint)

    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<w:pPr>' + `
             '<w:spacing w:before="320" w:after="320" w:line="360" w:lineRule="auto"/>' + `
             '<w:rPr>' + `
               '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
               '<w:sz w:val="24"/>' + `
               '<w:szCs w:val="24"/>' + `
               '<w:u w:val="single"/>' + `
             '</w:rPr>' + `
           '</w:pPr>' + `
           '</w:p>'

    $target.InsertXML($xml)
} else {
    Write-Output "Could not locate insertion point; document left unchanged."
}
